$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.689.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.42%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.166.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.99%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'529.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.93%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'135.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.35%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.16%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.166.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.447"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.23%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -6.99%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -7.22%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -7.28%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.705.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.05%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.129"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.42%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'25.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.160.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.91%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'58.613.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.64%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.0000153"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -6.32%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -4.67%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.38%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -6.50%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'345.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -7.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.09%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.512"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.77%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'67.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -7.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.288.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.97%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.173"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.15%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0952"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -6.94%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.88%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.60%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'Fetch.AI"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'1.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.11%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'PancakeSwap"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'1.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.97%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.66%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'21.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.26%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.75%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'159.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.96%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.14%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -8.92%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0690"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.47%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.195.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'40.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.64%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'23.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.45%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.699"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -6.55%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.68%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'3.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.10%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -6.09%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.296.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.84%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'6.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -3.87%  "
$ws.Range("E51").Style = "Normal"
